$d = $word.ActiveDocument

# --- Edit 1: update the subtitle line under the name ---
$d.Content.Find.Execute(
    "Data & BI Analyst | AI & Machine Learning | SQL, Python, R | ETL Pipelines, KPI Dashboards & Automation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data & BI Analyst | Healthcare & Fintech | AI & Analytics Engineering", 2) | Out-Null

# --- Edit 2: replace the trailing empty "List Paragraph" paragraph (just before
# the document's final blank paragraph) with the new "Key Portfolio Projects" section ---
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Style.NameLocal -eq "List Paragraph" -and $p.Range.Text.Trim() -eq "") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the empty List Paragraph placeholder to replace."
}

$target = $paras.Item($targetIndex).Range

$xmlFragment = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="121"/><w:rPr><w:color w:val="4471C4"/><w:spacing w:val="-2"/><w:w w:val="110"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="4471C4"/><w:spacing w:val="-2"/><w:w w:val="110"/></w:rPr><w:t>Key Portfolio Projects</w:t></w:r></w:p><w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>Fraud Detection Dashboard API | Personal Project | 2024</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve">Built production-ready </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>FastAPI</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> backend with ML monitoring, network analysis, and explainable AI capabilities; processes 5M+ transactions with sub-second response times, features 18 interactive analytics tiles, and enables real-time fraud ring detection and money mule account identification.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Demonstrated 3,700% ROI potential through comprehensive model monitoring, drift detection, and 24-hour predictive risk scoring for financial institutions.</w:t></w:r></w:p><w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>Bank Marketing Analytics Dashboard | Personal Project | 2024</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve">Developed full-stack ML analytics platform (React/Flask + </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Streamlit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>) featuring 6 ML algorithms for subscription prediction, achieving 90.5%+ accuracy, K-Means customer segmentation, and contact optimisation analysis.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Analysed 41,188+ customers to enable pre-call customer scoring, 15-20% conversion improvement through segmentation, and optimal contact frequency strategies for financial institutions.</w:t></w:r></w:p><w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>NIHR Research Intelligence Dashboard | University of Southampton | 2024</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve">Architected comprehensive BI platform (4,492 lines, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Streamlit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>/Python) with automated ETL, analysing £171.6M research portfolio across 11 research programmes and 314 constituencies.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:kern w:val="0"/><w:lang w:eastAsia="zh-CN"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> Demonstrated 122:1 ROI potential and £116K-£172K projected annual savings through 70% faster reporting, parliament-ready visualisations, and data quality scoring (73.6%).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xmlFragment)
